$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.635.31'
$ws.Range("E2").Value = '  +1.29%  '
$ws.Range("D3").Value = '3.449.59'
$ws.Range("E3").Value = '  +2.26%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.27%  '
$ws.Range("D7").Value = '3.450.42'
$ws.Range("E7").Value = '  +2.33%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.475'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.127'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.12%  '
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("D13").Value = '4.037.62'
$ws.Range("E13").Value = '  +2.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.92'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.48%  '
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("E16").Value = '  +2.19%  '
$ws.Range("D17").Value = '3.447.39'
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("D18").Value = '61.731.84'
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("E19").Value = '  +8.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '390.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.26%  '
$ws.Range("E23").Value = '  +2.69%  '
$ws.Range("D24").Value = '3.588.66'
$ws.Range("E24").Value = '  +1.96%  '
$ws.Range("E25").Value = '  +2.92%  '
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("E27").Value = '  +0.75%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("E29").Value = '  +2.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = '  -13.27%  '
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.74%  '
$ws.Range("E34").Value = '  +1.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '24.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.03'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.39%  '
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("E39").Value = '  +1.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '166.63'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0787'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.96'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +11.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.793'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.52'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.82%  '
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.31'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.70%  '
$ws.Range("E47").Value = '  +0.96%  '
$ws.Range("D48").Value = '2.609.25'
$ws.Range("E48").Value = '  +6.07%  '
$ws.Range("E49").Value = '  -1.86%  '
$ws.Range("E50").Value = '  +2.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.31%  '
